# Add the 2020 figures (a new column N) to the "Financial indicators of
# insurance companies" table on the active sheet, matching the formatting
# already used by the 2019 column (M), and move the active selection to N6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- N3: year header "2020" (same format as M3) ---
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial($xlPasteFormats)
$ws.Range("N3").Value = 2020

# --- N4: number of reporting insurance companies (same format as M4) ---
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial($xlPasteFormats)
$ws.Range("N4").Value = 15

# --- N5: insurance premiums, mln soms (same format as M5) ---
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial($xlPasteFormats)
$ws.Range("N5").Value = 1308.3

# Move the active selection to N6, just below the newly added column.
$ws.Range("N6").Select()
